$wb = $excel.ActiveWorkbook

# --- Linear sheet ---
$wsLinear = $wb.Worksheets.Item("Linear")
$wsLinear.Range("B2").Value = 3.028490980480218
$wsLinear.Range("B3").Value = 0.1078021605665619
$wsLinear.Range("B4").Value = 35985.49555355266
$wsLinear.Range("B5").Value = "[1.0, 0.17010076790456002, 0.010706161581064667, 0.03717522595810434, 0.017287988258759334, -0.052919377798628646, 0.10802663114565386, 0.254868183895977, 0.08788936227154939, -0.04302368268070455, -0.0033739226686931504, 0.0037877225469762996, -0.05003684191673645, 0.09208468422852623, 0.23484164540068314, 0.05328825968606437, -0.0710543490897814, -0.016102352282465236, -0.01571954214061154, -0.03592944028630505]"

# --- NonLinear sheet ---
$wsNonLinear = $wb.Worksheets.Item("NonLinear")
$wsNonLinear.Range("B4").Value = -4.863491300520318
$wsNonLinear.Range("B5").Value = 0.02830574927192976
$wsNonLinear.Range("B6").Value = 35124.73483759602
$wsNonLinear.Range("B7").Value = 5.497741303112915
$wsNonLinear.Range("B8").Value = 0.1063036189697222
$wsNonLinear.Range("B9").Value = 36791.43935766975
$wsNonLinear.Range("B10").Value = "[1.0, 0.1677044748924979, 0.008683634589049218, 0.03442552396094724, 0.015310570794049334, -0.05325621351830622, 0.10735261046526622, 0.2514731834790395, 0.08727851971451552, -0.04451823269005174, -0.004002530847665953, 0.004521748815596559, -0.0499305570245723, 0.09126907797189042, 0.23216408560972449, 0.05312699460311114, -0.07065742628160213, -0.017207026597119483, -0.015099095334449716, -0.034541823793790845]"
